# Applies the scheduled-runner market-data refresh to the Leve profit tables.
# For each sheet, update currentAveragePrice/NQ/HQ, LevePrice NQ/HQ and
# LeveProfit NQ/HQ (columns H:N) for the rows whose source data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 85
$ws.Range("I4").Value = 85
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 85
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 29
$ws.Range("N4").ClearContents()
$ws.Range("H57").Value = 34139.332
$ws.Range("J57").Value = 34139.332
$ws.Range("L57").Value = 102417.996
$ws.Range("N57").Value = -103415.996
$ws.Range("H86").Value = 5320
$ws.Range("I86").Value = 4640
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 4640
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -3517
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 5320
$ws.Range("I89").Value = 4640
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 23200
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -17584
$ws.Range("N89").Value = -41232
$ws.Range("H106").Value = 1424.3
$ws.Range("I106").Value = 693.6667
$ws.Range("K106").Value = 693.6667
$ws.Range("M106").Value = -62.66669999999999
$ws.Range("H132").Value = 2990.6428
$ws.Range("I132").Value = 2937.658
$ws.Range("J132").Value = 3494
$ws.Range("K132").Value = 8812.974
$ws.Range("L132").Value = 10482
$ws.Range("M132").Value = -6282.974
$ws.Range("N132").Value = -15542
$ws.Range("H137").Value = 11620
$ws.Range("I137").Value = 1791.8334
$ws.Range("K137").Value = 5375.5002
$ws.Range("M137").Value = -2825.5002
$ws.Range("H138").Value = 8716
$ws.Range("J138").Value = 11945.094
$ws.Range("L138").Value = 35835.282
$ws.Range("N138").Value = -46115.282
$ws.Range("H141").Value = 5792.9165
$ws.Range("I141").Value = 6589.4443
$ws.Range("K141").Value = 19768.3329
$ws.Range("M141").Value = -14588.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14289100
$ws.Range("I32").Value = 16394965
$ws.Range("J32").Value = 16018
$ws.Range("K32").Value = 16394965
$ws.Range("L32").Value = 16018
$ws.Range("M32").Value = -16394678
$ws.Range("N32").Value = -16592
$ws.Range("H97").Value = 1788.0834
$ws.Range("I97").Value = 1757.1305
$ws.Range("J97").Value = 2500
$ws.Range("K97").Value = 1757.1305
$ws.Range("L97").Value = 2500
$ws.Range("M97").Value = -1261.1305
$ws.Range("N97").Value = -3492

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2662.6858
$ws.Range("I20").Value = 2106
$ws.Range("J20").Value = 3323.75
$ws.Range("K20").Value = 2106
$ws.Range("L20").Value = 3323.75
$ws.Range("M20").Value = -1859
$ws.Range("N20").Value = -3817.75
$ws.Range("H22").Value = 253.57143
$ws.Range("I22").Value = 253.57143
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 253.57143
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -80.57142999999999
$ws.Range("N22").ClearContents()
$ws.Range("H140").Value = 53806.332
$ws.Range("J140").Value = 53806.332
$ws.Range("L140").Value = 53806.332
$ws.Range("N140").Value = -64166.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 916.5714
$ws.Range("I16").Value = 923.0909
$ws.Range("K16").Value = 923.0909
$ws.Range("M16").Value = -636.0909
$ws.Range("H31").Value = 46299564
$ws.Range("I31").Value = 3638.75
$ws.Range("J31").Value = 65792588
$ws.Range("K31").Value = 3638.75
$ws.Range("L31").Value = 65792588
$ws.Range("M31").Value = -3343.75
$ws.Range("N31").Value = -65793178
$ws.Range("H34").Value = 46299564
$ws.Range("I34").Value = 3638.75
$ws.Range("J34").Value = 65792588
$ws.Range("K34").Value = 3638.75
$ws.Range("L34").Value = 65792588
$ws.Range("M34").Value = -3436.75
$ws.Range("N34").Value = -65792992
$ws.Range("H113").Value = 916.5714
$ws.Range("I113").Value = 923.0909
$ws.Range("K113").Value = 923.0909
$ws.Range("M113").Value = 1246.9091

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4244.7
$ws.Range("J3").Value = 2100
$ws.Range("L3").Value = 6300
$ws.Range("N3").Value = -6524
$ws.Range("H68").Value = 5549.225
$ws.Range("I68").Value = 8249.25
$ws.Range("J68").Value = 5249.222
$ws.Range("K68").Value = 24747.75
$ws.Range("L68").Value = 15747.666
$ws.Range("M68").Value = -23936.75
$ws.Range("N68").Value = -17369.666
$ws.Range("H71").Value = 5549.225
$ws.Range("I71").Value = 8249.25
$ws.Range("J71").Value = 5249.222
$ws.Range("K71").Value = 74243.25
$ws.Range("L71").Value = 47242.998
$ws.Range("M71").Value = -70187.25
$ws.Range("N71").Value = -55354.998
$ws.Range("H132").Value = 5559976
$ws.Range("I132").Value = 1649.75
$ws.Range("J132").Value = 8339139
$ws.Range("K132").Value = 14847.75
$ws.Range("L132").Value = 75052251
$ws.Range("M132").Value = -12317.75
$ws.Range("N132").Value = -75057311

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8199.875
$ws.Range("I80").Value = 16750
$ws.Range("J80").Value = 5349.8335
$ws.Range("K80").Value = 16750
$ws.Range("L80").Value = 5349.8335
$ws.Range("M80").Value = -15752
$ws.Range("N80").Value = -7345.8335
$ws.Range("H83").Value = 8199.875
$ws.Range("I83").Value = 16750
$ws.Range("J83").Value = 5349.8335
$ws.Range("K83").Value = 83750
$ws.Range("L83").Value = 26749.1675
$ws.Range("M83").Value = -78758
$ws.Range("N83").Value = -36733.1675

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1447.4
$ws.Range("I46").Value = 1018.25
$ws.Range("J46").Value = 3164
$ws.Range("K46").Value = 1018.25
$ws.Range("L46").Value = 3164
$ws.Range("M46").Value = -830.25
$ws.Range("N46").Value = -3540
$ws.Range("H122").Value = 2607446
$ws.Range("I122").Value = 2921.2896
$ws.Range("J122").Value = 12504640
$ws.Range("K122").Value = 8763.8688
$ws.Range("L122").Value = 37513920
$ws.Range("M122").Value = -6313.8688
$ws.Range("N122").Value = -37518820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 24
$ws.Range("I25").Value = 24
$ws.Range("K25").Value = 24
$ws.Range("M25").Value = 269
$ws.Range("H136").Value = 1321.7446
$ws.Range("I136").Value = 1345.3513
$ws.Range("J136").Value = 1234.4
$ws.Range("K136").Value = 4036.0539
$ws.Range("L136").Value = 3703.2
$ws.Range("M136").Value = -1486.0539
$ws.Range("N136").Value = -8803.200000000001

